# Update "想去人数" (want-to-go count) figures across the workbook sheets.
# These figures appear once in their "home" sheet (展览 / 演出 / 本地生活)
# and are duplicated again inside the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 832
$ws1.Range("F8").Value  = 219
$ws1.Range("F12").Value = 653
$ws1.Range("F15").Value = 349
$ws1.Range("F16").Value = 3045
$ws1.Range("F19").Value = 52

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 92

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5319
$ws3.Range("F4").Value = 250

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5319
$ws4.Range("F6").Value  = 250
$ws4.Range("F13").Value = 832
$ws4.Range("F19").Value = 219
$ws4.Range("F26").Value = 653
$ws4.Range("F28").Value = 92
$ws4.Range("F30").Value = 349
$ws4.Range("F31").Value = 3045
$ws4.Range("F35").Value = 52
